$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 11 - Bourgogne-Franche-Comte / SARL (54)
Set-TextValue "C11" "564"
Set-TextValue "D11" "493"
Set-TextValue "E11" "4545867.99"

# Row 12 - Bourgogne-Franche-Comte / SAS (57)
Set-TextValue "C12" "279"
Set-TextValue "E12" "2298317.49"

# Row 27 - Corse / Entrepreneur individuel (10)
Set-TextValue "C27" "153"
Set-TextValue "D27" "148"
Set-TextValue "E27" "403130.37"

# Row 28 - Corse / SARL (54)
Set-TextValue "C28" "174"
Set-TextValue "E28" "937067.73"

# Row 29 - Corse / SAS (57)
Set-TextValue "C29" "150"
Set-TextValue "E29" "703789.60"

# Row 65 - La Reunion / SAS (57)
Set-TextValue "C65" "61"
Set-TextValue "E65" "470050.00"

# Row 81 - Nouvelle-Aquitaine / SARL (54)
Set-TextValue "C81" "1308"
Set-TextValue "D81" "1111"
Set-TextValue "E81" "10868143.05"
